$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Interventions coverages")

# Header row: rename existing "pre-2016" header to "baseline coverage" and add two new
# headers for the saturation coverage and unit cost columns.
$ws.Range("B1").Value = "baseline coverage"
$ws.Range("C1").Value = "saturation coverage"
$ws.Range("D1").Value = "unit cost"

# Saturation coverage (column C) and unit cost (column D) data for each intervention.
$ws.Range("C2").Value = 0.85
$ws.Range("D2").Value = 60

$ws.Range("C3").Value = 0.85
$ws.Range("D3").Value = 300

$ws.Range("C4").Value = 0.85
$ws.Range("D4").Value = 300

$ws.Range("C5").Value = 0.85
$ws.Range("D5").Value = 100

$ws.Range("C6").Value = 0.85
$ws.Range("D6").Value = 20

$ws.Range("C7").Value = 0.85
$ws.Range("D7").Value = 80

$ws.Range("C8").Value = 0.85
$ws.Range("D8").Value = 90

# Widen the columns to fit the new headers/content.
$ws.Columns.Item(2).ColumnWidth = 20.71
$ws.Columns.Item(3).ColumnWidth = 21.43
$ws.Columns.Item(4).ColumnWidth = 21.43

# "Interventions maternal" sheet: fill in the affected-fraction row for IPTp (row 3),
# which previously had all zeros.
$ws2 = $wb.Worksheets.Item("Interventions maternal")
$ws2.Range("C3").Value = 0.1
$ws2.Range("D3").Value = 0.1
$ws2.Range("E3").Value = 0.1
$ws2.Range("F3").Value = 0.1
